# EMEP_NFR14_scaling_mapping_SO2.xlsx
# "Updated S scaling to improve match to inventory data in 1980 and before."
#
# Changes applied on the "method" sheet:
#   * The explanatory comment text shown in the header (G1) is extended with a
#     note about using 1990 as the scaling year for Eastern Europe + esp
#     Industry/Power.
#   * esp (Spain) Industry/PublicPower rows: pre-extension year 1980 -> 1990,
#     pre-extension method linear_1 -> constant.
#   * lva/ltu/geo/est Industry & PublicPower rows: interpolation method
#     linear -> constant.
#   * svn/mkd/hrv Industry & PublicPower rows: pre-extension method
#     linear_1 -> constant.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("method")

# --- Header comment (column G, row 1) ---------------------------------------
$ws.Range("G1").Value = "(Better match to inventory if don't use linear to 1 for hun, bgr, gbr ind/power. Use 1990 year for EE + esp for Industry/power."

# --- esp (Spain): pre_ext_year 1980 -> 1990, pre_ext_method linear_1 -> constant
foreach ($r in 28, 29) {
    $ws.Range("C$r").Value = 1990
    $ws.Range("D$r").Value = "constant"
}

# --- interp_method (column E): linear -> constant
foreach ($r in 76, 80, 81, 84, 85, 88, 89) {
    $ws.Range("E$r").Value = "constant"
}

# --- pre_ext_method (column D): linear_1 -> constant
foreach ($r in 111, 112, 115, 116, 119, 120) {
    $ws.Range("D$r").Value = "constant"
}

# --- Leave the sheet's selection on the last-edited cell, as in the source edit
$ws.Activate() | Out-Null
$ws.Range("A112").Select() | Out-Null
